$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B12").Value = 0.0305653999967034
$ws.Range("C12").Value = 0.00049919998855329996

$ws.Activate()
$ws.Range("E12").Select()
